# This workbook edit re-orders the observation rows 3-14 on the "Artfynd"
# sheet (row 9 keeps its original content and position). For every other
# row we overwrite the observation-specific columns (Id, Taxonsorterings-
# ordning, Rödlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor,
# Aktivitet, Ost, Nord, Substrat-beskrivning) with the values that used to
# belong to a different row, per the source material. All other columns
# (dates, location names, reporter, etc.) are identical across these rows
# already, so they do not need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current ("before") values up front, so that writing the new
# value into one row never clobbers data that another row still needs to read.
$snap = @{}
$snap['A3'] = 111473774
$snap['B3'] = 89405
$snap['D3'] = 'NT'
$snap['E3'] = 1202
$snap['F3'] = 'Ullticka'
$snap['G3'] = 'Phellinidium ferrugineofuscum'
$snap['H3'] = '(P.Karst.) Fiasson & Niemelä'
$snap['Q3'] = 703999.5190368021
$snap['R3'] = 6572850.823973293
$snap['AO3'] = 'granlåga'
$snap['A4'] = 111473783
$snap['B4'] = 89686
$snap['D4'] = 'NT'
$snap['E4'] = 658
$snap['F4'] = 'Rosenticka'
$snap['G4'] = 'Rhodofomes roseus'
$snap['H4'] = '(Alb. & Schwein.) Kotl. & Pouzar'
$snap['Q4'] = 703998.3853129407
$snap['R4'] = 6572852.813158008
$snap['AO4'] = 'granlåga'
$snap['A5'] = 111473792
$snap['B5'] = 5113
$snap['D5'] = 'LC'
$snap['E5'] = 100526
$snap['F5'] = 'Bronshjon'
$snap['G5'] = 'Callidium coriaceum'
$snap['H5'] = 'Paykull, 1800'
$snap['M5'] = 'äldre gnagspår'
$snap['Q5'] = 703965.55072247
$snap['R5'] = 6572785.445717536
$snap['AO5'] = 'torrgran'
$snap['A6'] = 111473779
$snap['B6'] = 89425
$snap['D6'] = 'NT'
$snap['E6'] = 5442
$snap['F6'] = 'Tallticka'
$snap['G6'] = 'Porodaedalea pini'
$snap['H6'] = '(Brot.) Murrill'
$snap['Q6'] = 704193.4830821306
$snap['R6'] = 6572948.378178579
$snap['AO6'] = 'gammeltall'
$snap['A7'] = 111473777
$snap['B7'] = 89425
$snap['D7'] = 'NT'
$snap['E7'] = 5442
$snap['F7'] = 'Tallticka'
$snap['G7'] = 'Porodaedalea pini'
$snap['H7'] = '(Brot.) Murrill'
$snap['Q7'] = 704301.1177162804
$snap['R7'] = 6573209.392206083
$snap['AO7'] = 'gammeltall'
$snap['A8'] = 111473784
$snap['B8'] = 73634
$snap['D8'] = 'LC'
$snap['E8'] = 6426
$snap['F8'] = 'Kattfotslav'
$snap['G8'] = 'Felipes leucopellaeus'
$snap['H8'] = '(Ach.) Frisch & G.Thor'
$snap['Q8'] = 704135.470341172
$snap['R8'] = 6572843.267234835
$snap['AO8'] = 'äldre gran'
$snap['A10'] = 111473791
$snap['B10'] = 93289
$snap['D10'] = 'LC'
$snap['E10'] = 2170
$snap['F10'] = 'Flagellkvastmossa'
$snap['G10'] = 'Dicranum flagellare'
$snap['H10'] = 'Hedw.'
$snap['Q10'] = 704004.9502936595
$snap['R10'] = 6572835.740028554
$snap['AO10'] = 'låga av tall'
$snap['A11'] = 111473782
$snap['B11'] = 89183
$snap['D11'] = 'LC'
$snap['E11'] = 3215
$snap['F11'] = 'Rödgul trumpetsvamp'
$snap['G11'] = 'Craterellus lutescens'
$snap['H11'] = '(Fr.) Fr.'
$snap['Q11'] = 704171.5165585374
$snap['R11'] = 6572850.843097115
$snap['A12'] = 111473793
$snap['B12'] = 93388
$snap['D12'] = 'LC'
$snap['E12'] = 2180
$snap['F12'] = 'Blåmossa'
$snap['G12'] = 'Leucobryum glaucum'
$snap['H12'] = '(Hedw.) Ångstr.'
$snap['Q12'] = 703959.3331032015
$snap['R12'] = 6572805.612961343
$snap['A13'] = 111473775
$snap['B13'] = 89405
$snap['D13'] = 'NT'
$snap['E13'] = 1202
$snap['F13'] = 'Ullticka'
$snap['G13'] = 'Phellinidium ferrugineofuscum'
$snap['H13'] = '(P.Karst.) Fiasson & Niemelä'
$snap['Q13'] = 703969.3444121893
$snap['R13'] = 6572791.287347207
$snap['AO13'] = 'granlåga'
$snap['A14'] = 111473773
$snap['B14'] = 89405
$snap['D14'] = 'NT'
$snap['E14'] = 1202
$snap['F14'] = 'Ullticka'
$snap['G14'] = 'Phellinidium ferrugineofuscum'
$snap['H14'] = '(P.Karst.) Fiasson & Niemelä'
$snap['Q14'] = 704016.0051346947
$snap['R14'] = 6572801.994589122
$snap['AO14'] = 'granlåga'

# Apply the permutation: each destination row receives the snapshotted
# values from its designated source row.

# Row 3 <= old row 10
$ws.Range("A3").Value2 = $snap['A10']
$ws.Range("B3").Value2 = $snap['B10']
$ws.Range("D3").Value2 = $snap['D10']
$ws.Range("E3").Value2 = $snap['E10']
$ws.Range("F3").Value2 = $snap['F10']
$ws.Range("G3").Value2 = $snap['G10']
$ws.Range("H3").Value2 = $snap['H10']
$ws.Range("M3").ClearContents()
$ws.Range("Q3").Value2 = $snap['Q10']
$ws.Range("R3").Value2 = $snap['R10']
$ws.Range("AO3").Value2 = $snap['AO10']

# Row 4 <= old row 5
$ws.Range("A4").Value2 = $snap['A5']
$ws.Range("B4").Value2 = $snap['B5']
$ws.Range("D4").Value2 = $snap['D5']
$ws.Range("E4").Value2 = $snap['E5']
$ws.Range("F4").Value2 = $snap['F5']
$ws.Range("G4").Value2 = $snap['G5']
$ws.Range("H4").Value2 = $snap['H5']
$ws.Range("M4").Value2 = $snap['M5']
$ws.Range("Q4").Value2 = $snap['Q5']
$ws.Range("R4").Value2 = $snap['R5']
$ws.Range("AO4").Value2 = $snap['AO5']

# Row 5 <= old row 11
$ws.Range("A5").Value2 = $snap['A11']
$ws.Range("B5").Value2 = $snap['B11']
$ws.Range("D5").Value2 = $snap['D11']
$ws.Range("E5").Value2 = $snap['E11']
$ws.Range("F5").Value2 = $snap['F11']
$ws.Range("G5").Value2 = $snap['G11']
$ws.Range("H5").Value2 = $snap['H11']
$ws.Range("M5").ClearContents()
$ws.Range("Q5").Value2 = $snap['Q11']
$ws.Range("R5").Value2 = $snap['R11']
$ws.Range("AO5").ClearContents()

# Row 6 <= old row 13
$ws.Range("A6").Value2 = $snap['A13']
$ws.Range("B6").Value2 = $snap['B13']
$ws.Range("D6").Value2 = $snap['D13']
$ws.Range("E6").Value2 = $snap['E13']
$ws.Range("F6").Value2 = $snap['F13']
$ws.Range("G6").Value2 = $snap['G13']
$ws.Range("H6").Value2 = $snap['H13']
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value2 = $snap['Q13']
$ws.Range("R6").Value2 = $snap['R13']
$ws.Range("AO6").Value2 = $snap['AO13']

# Row 7 <= old row 6
$ws.Range("A7").Value2 = $snap['A6']
$ws.Range("B7").Value2 = $snap['B6']
$ws.Range("D7").Value2 = $snap['D6']
$ws.Range("E7").Value2 = $snap['E6']
$ws.Range("F7").Value2 = $snap['F6']
$ws.Range("G7").Value2 = $snap['G6']
$ws.Range("H7").Value2 = $snap['H6']
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value2 = $snap['Q6']
$ws.Range("R7").Value2 = $snap['R6']
$ws.Range("AO7").Value2 = $snap['AO6']

# Row 8 <= old row 12
$ws.Range("A8").Value2 = $snap['A12']
$ws.Range("B8").Value2 = $snap['B12']
$ws.Range("D8").Value2 = $snap['D12']
$ws.Range("E8").Value2 = $snap['E12']
$ws.Range("F8").Value2 = $snap['F12']
$ws.Range("G8").Value2 = $snap['G12']
$ws.Range("H8").Value2 = $snap['H12']
$ws.Range("M8").ClearContents()
$ws.Range("Q8").Value2 = $snap['Q12']
$ws.Range("R8").Value2 = $snap['R12']
$ws.Range("AO8").ClearContents()

# Row 10 <= old row 14
$ws.Range("A10").Value2 = $snap['A14']
$ws.Range("B10").Value2 = $snap['B14']
$ws.Range("D10").Value2 = $snap['D14']
$ws.Range("E10").Value2 = $snap['E14']
$ws.Range("F10").Value2 = $snap['F14']
$ws.Range("G10").Value2 = $snap['G14']
$ws.Range("H10").Value2 = $snap['H14']
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value2 = $snap['Q14']
$ws.Range("R10").Value2 = $snap['R14']
$ws.Range("AO10").Value2 = $snap['AO14']

# Row 11 <= old row 8
$ws.Range("A11").Value2 = $snap['A8']
$ws.Range("B11").Value2 = $snap['B8']
$ws.Range("D11").Value2 = $snap['D8']
$ws.Range("E11").Value2 = $snap['E8']
$ws.Range("F11").Value2 = $snap['F8']
$ws.Range("G11").Value2 = $snap['G8']
$ws.Range("H11").Value2 = $snap['H8']
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value2 = $snap['Q8']
$ws.Range("R11").Value2 = $snap['R8']
$ws.Range("AO11").Value2 = $snap['AO8']

# Row 12 <= old row 3
$ws.Range("A12").Value2 = $snap['A3']
$ws.Range("B12").Value2 = $snap['B3']
$ws.Range("D12").Value2 = $snap['D3']
$ws.Range("E12").Value2 = $snap['E3']
$ws.Range("F12").Value2 = $snap['F3']
$ws.Range("G12").Value2 = $snap['G3']
$ws.Range("H12").Value2 = $snap['H3']
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value2 = $snap['Q3']
$ws.Range("R12").Value2 = $snap['R3']
$ws.Range("AO12").Value2 = $snap['AO3']

# Row 13 <= old row 7
$ws.Range("A13").Value2 = $snap['A7']
$ws.Range("B13").Value2 = $snap['B7']
$ws.Range("D13").Value2 = $snap['D7']
$ws.Range("E13").Value2 = $snap['E7']
$ws.Range("F13").Value2 = $snap['F7']
$ws.Range("G13").Value2 = $snap['G7']
$ws.Range("H13").Value2 = $snap['H7']
$ws.Range("M13").ClearContents()
$ws.Range("Q13").Value2 = $snap['Q7']
$ws.Range("R13").Value2 = $snap['R7']
$ws.Range("AO13").Value2 = $snap['AO7']

# Row 14 <= old row 4
$ws.Range("A14").Value2 = $snap['A4']
$ws.Range("B14").Value2 = $snap['B4']
$ws.Range("D14").Value2 = $snap['D4']
$ws.Range("E14").Value2 = $snap['E4']
$ws.Range("F14").Value2 = $snap['F4']
$ws.Range("G14").Value2 = $snap['G4']
$ws.Range("H14").Value2 = $snap['H4']
$ws.Range("M14").ClearContents()
$ws.Range("Q14").Value2 = $snap['Q4']
$ws.Range("R14").Value2 = $snap['R4']
$ws.Range("AO14").Value2 = $snap['AO4']
